$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 1.25
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 2.72
$ws.Range("Q2").Value = 1.52
$ws.Range("R2").Value = 1.7
$ws.Range("S2").Value = 2.26
$ws.Range("T2").Value = 1.52
$ws.Range("U2").Value = 2.74
$ws.Range("X2").Value = 27
$ws.Range("Y2").Value = 23
$ws.Range("AB2").Value = 16
$ws.Range("AC2").Value = 10.5
$ws.Range("AF2").Value = 17
$ws.Range("AH2").Value = 14.5
$ws.Range("AJ2").Value = 27
$ws.Range("AM2").Value = 55
$ws.Range("AN2").Value = 9.199999999999999
$ws.Range("AO2").Value = 21
# Row 3
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2.14
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3.85
$ws.Range("K3").Value = 4.3
$ws.Range("N3").Value = 4.8
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 2.3
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 1.53
$ws.Range("S3").Value = 2.62
$ws.Range("T3").Value = 1.61
$ws.Range("U3").Value = 2.42
$ws.Range("W3").Value = 1.87
$ws.Range("X3").Value = 26
$ws.Range("Y3").Value = 22
$ws.Range("Z3").Value = 36
$ws.Range("AB3").Value = 15.5
$ws.Range("AC3").Value = 9.6
$ws.Range("AD3").Value = 19
$ws.Range("AE3").Value = 46
$ws.Range("AF3").Value = 18.5
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 18
$ws.Range("AI3").Value = 50
$ws.Range("AJ3").Value = 30
$ws.Range("AK3").Value = 23
$ws.Range("AL3").Value = 36
$ws.Range("AN3").Value = 13.5
$ws.Range("AO3").Value = 36
# Row 4
$ws.Range("F4").Value = 5.4
$ws.Range("G4").Value = 6.4
$ws.Range("H4").Value = 1.66
$ws.Range("I4").Value = 1.75
$ws.Range("L4").Value = 1.32
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 3.8
$ws.Range("P4").Value = 1.95
$ws.Range("Q4").Value = 1.73
$ws.Range("R4").Value = 1.37
$ws.Range("S4").Value = 3.2
$ws.Range("T4").Value = 1.84
$ws.Range("V4").Value = 2.32
$ws.Range("W4").Value = 1.19
$ws.Range("X4").Value = 19.5
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 11
$ws.Range("AA4").Value = 20
$ws.Range("AB4").Value = 21
$ws.Range("AC4").Value = 11
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 21
$ws.Range("AF4").Value = 60
$ws.Range("AG4").Value = 24
$ws.Range("AH4").Value = 22
$ws.Range("AJ4").Value = 190
$ws.Range("AK4").Value = 100
$ws.Range("AM4").Value = 140
$ws.Range("AO4").Value = 11
# Row 6
$ws.Range("F6").Value = 3
$ws.Range("H6").Value = 2.68
$ws.Range("I6").Value = 2.86
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.25
$ws.Range("L6").Value = 1.6
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 2.5
$ws.Range("O6").Value = 1.59
$ws.Range("P6").Value = 1.51
$ws.Range("Q6").Value = 2.74
$ws.Range("R6").Value = 1.19
$ws.Range("S6").Value = 5.8
$ws.Range("T6").Value = 2.18
$ws.Range("U6").Value = 1.73
$ws.Range("W6").Value = 1.45
$ws.Range("X6").Value = 10
$ws.Range("Y6").Value = 7.8
$ws.Range("Z6").Value = 980
$ws.Range("AA6").Value = 980
$ws.Range("AB6").Value = 8.6
$ws.Range("AC6").Value = 8.4
$ws.Range("AD6").Value = 15.5
$ws.Range("AE6").Value = 980
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 29
$ws.Range("AI6").Value = 70
$ws.Range("AK6").Value = 55
$ws.Range("AL6").Value = 90
$ws.Range("AM6").Value = 220
$ws.Range("AN6").Value = 75
$ws.Range("AO6").Value = 60
# Row 7
$ws.Range("J7").Value = 1.03
$ws.Range("N7").Value = 1.3
$ws.Range("P7").Value = 1.3
$ws.Range("T7").Value = 1.03
$ws.Range("U7").Value = 1.03
# Row 8
$ws.Range("F8").Value = 1.45
$ws.Range("G8").Value = 1.47
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 12.5
$ws.Range("J8").Value = 4.4
$ws.Range("K8").Value = 5
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 3.35
$ws.Range("O8").Value = 1.34
$ws.Range("P8").Value = 1.81
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.3
$ws.Range("S8").Value = 3.6
$ws.Range("T8").Value = 2.3
$ws.Range("U8").Value = 1.63
$ws.Range("V8").Value = 1.09
$ws.Range("W8").Value = 3.1
$ws.Range("X8").Value = 14
$ws.Range("Y8").Value = 29
$ws.Range("Z8").Value = 110
$ws.Range("AA8").Value = 580
$ws.Range("AB8").Value = 7
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 44
$ws.Range("AE8").Value = 270
$ws.Range("AF8").Value = 7.6
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 260
$ws.Range("AI8").Value = 250
$ws.Range("AJ8").Value = 12
$ws.Range("AK8").Value = 18
$ws.Range("AL8").Value = 55
$ws.Range("AM8").Value = 320
$ws.Range("AN8").Value = 8.800000000000001
$ws.Range("AO8").Value = 480
# Row 9
$ws.Range("F9").Value = 2.98
$ws.Range("G9").Value = 3.2
$ws.Range("H9").Value = 2.7
$ws.Range("I9").Value = 2.86
$ws.Range("J9").Value = 3.1
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 2.66
$ws.Range("O9").Value = 1.54
$ws.Range("P9").Value = 1.54
$ws.Range("Q9").Value = 2.64
$ws.Range("R9").Value = 1.19
$ws.Range("S9").Value = 5.4
$ws.Range("T9").Value = 2.06
$ws.Range("U9").Value = 1.77
$ws.Range("V9").Value = 1.53
$ws.Range("W9").Value = 1.46
$ws.Range("Y9").Value = 8.6
$ws.Range("AB9").Value = 9.199999999999999
$ws.Range("AC9").Value = 7.4
$ws.Range("AF9").Value = 18.5
$ws.Range("AG9").Value = 15
$ws.Range("AK9").Value = 50
$ws.Range("AO9").Value = 48
# Row 10
$ws.Range("F10").Value = 1.35
$ws.Range("G10").Value = 1.4
$ws.Range("H10").Value = 11
$ws.Range("I10").Value = 14.5
$ws.Range("J10").Value = 5.2
$ws.Range("K10").Value = 5.4
$ws.Range("N10").Value = 3.35
$ws.Range("P10").Value = 1.81
$ws.Range("Q10").Value = 2.06
$ws.Range("R10").Value = 1.29
$ws.Range("S10").Value = 3.75
$ws.Range("T10").Value = 2.48
$ws.Range("U10").Value = 1.56
$ws.Range("W10").Value = 3.5
$ws.Range("AB10").Value = 6.6
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 55
$ws.Range("AK10").Value = 18
$ws.Range("AL10").Value = 70
$ws.Range("AN10").Value = 8.199999999999999
$ws.Range("AO10").Value = 880
# Row 11
$ws.Range("J11").Value = 1.03
$ws.Range("N11").Value = 1.36
$ws.Range("P11").Value = 1.36
$ws.Range("R11").Value = 1.18
$ws.Range("T11").Value = 1.03
$ws.Range("U11").Value = 1.03
